$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: 001 -> 002 (keep as text, avoid leaving a lingering explicit style)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("J2").ClearFormats()

# N2: report date
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = "2020-06-30 00:00:00"
$ws.Range("N2").ClearFormats()

# O2:AG2 numeric updates
$ws.Range("O2").Value = 1298832753.92
$ws.Range("P2").Value = 173050314.84
$ws.Range("Q2").Value = 635208993.71
$ws.Range("R2").Value = 28.3778662123
$ws.Range("S2").Value = 61657361.09
$ws.Range("T2").Value = -45.5213939836
$ws.Range("U2").Value = 124317361.89
$ws.Range("V2").Value = 118.0109340664
$ws.Range("W2").Value = 122632144.8
$ws.Range("X2").Value = 32258180.59
$ws.Range("Y2").Value = -4.1139466561

# Z2, AA2: cleared to blank (advance receivables figures no longer reported)
$ws.Range("Z2").Value = "'"
$ws.Range("Z2").ClearFormats()
$ws.Range("AA2").Value = "'"
$ws.Range("AA2").ClearFormats()

$ws.Range("AB2").Value = 1176200609.12
$ws.Range("AC2").Value = 7.951101486
$ws.Range("AD2").Value = 10.4843875668
$ws.Range("AE2").Value = 42.5749627737
$ws.Range("AF2").Value = 998.9568345439
$ws.Range("AG2").Value = 9.441719453899999
